$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.841.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9940"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6244"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9939"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07419"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2911"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.87"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07657"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.809.55"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.969"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6634"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.51"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009569"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.004"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.824.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "223.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9942"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.066"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9954"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.56"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1400"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.445"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.79"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.489"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.094"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.026"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05433"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.189"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.841"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7379"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.128"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.596"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.222.74"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.735"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01771"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.632"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8917"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9936"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.10"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.72"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000122"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5056"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4014"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.923"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.650"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "XinFinNetwork"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07167"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05777"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.46%  "
